# ---------------------------------------------------------------------------
# Add an "actual_size" column (measured particle size in weeks 1-2) between
# the existing "size" and "video" columns, fill in the selected-particle
# measurements, hide the rows of particles that were not selected ("usable"
# = 0) and turn on the AutoFilter that was used to pick the usable rows.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new column F ("actual_size"); this pushes the existing
#    "video" (F) and "usable" (G) columns one place to the right.
$ws.Columns("F:F").Insert()
$ws.Range("F1").Value = "actual_size"

# 2. Fill in the measured actual size for every particle that was marked
#    usable (column H after the insert). Particles that were not selected
#    are left blank in this column.
$ws.Range("F2").Value = 11.5
$ws.Range("F3").Value = 25.0
$ws.Range("F8").Value = 22.5
$ws.Range("F11").Value = 12.5
$ws.Range("F13").Value = 20.5
$ws.Range("F14").Value = 10.0
$ws.Range("F15").Value = 9.5
$ws.Range("F17").Value = 8.5
$ws.Range("F18").Value = 16.0
$ws.Range("F20").Value = 11.5
$ws.Range("F22").Value = 13.0
$ws.Range("F23").Value = 10.5
$ws.Range("F26").Value = 17.0
$ws.Range("F30").Value = 35.0
$ws.Range("F32").Value = 14.5

# 3. Widen the columns to fit the new data.
$ws.Columns("C:D").ColumnWidth = 11.16
$ws.Columns("F:F").ColumnWidth = 9.65

# 4. Hide the rows for particles that were not selected (usable = 0).
$ws.Rows("4:7").Hidden = $true
$ws.Rows("9:10").Hidden = $true
$ws.Rows(12).Hidden = $true
$ws.Rows(16).Hidden = $true
$ws.Rows(19).Hidden = $true
$ws.Rows(21).Hidden = $true
$ws.Rows("24:25").Hidden = $true
$ws.Rows("27:29").Hidden = $true
$ws.Rows(31).Hidden = $true
$ws.Rows(33).Hidden = $true

# 5. Turn on the AutoFilter over the whole table, filtered on the "usable"
#    column (now H) to show only the selected (1) rows.
$usableValues = @("1")
$ws.Range("A1:H33").AutoFilter(8, $usableValues, 7)

# 6. Record the filter range as the sheet's (hidden) _FilterDatabase name,
#    as Excel does when AutoFilter is applied.
$fd = $ws.Names.Add("_xlnm._FilterDatabase", "='7_particle_selection'!`$A`$1:`$H`$33")
$fd.Visible = $false
